# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the per-language sheets to reflect the newly generated
# handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 04:36:24"
$wsZhCn.Range("H2").Value = "2016-03-18 04:36:41"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 04:36:30"
$wsDeDe.Range("H2").Value = "2016-03-18 04:36:47"
